$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.650.70'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.035.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.55%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.96'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.606'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.63'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.28%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.375'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0832'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.337.07'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.38'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.94%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.10'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.50'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.771'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.048.06'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.638.36'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.96'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.26'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0822'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.79'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.47%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.40'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.72'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.76'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.32%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.24'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +9.61%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0608'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.40%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.51'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.34'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.40'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.97%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.96'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.530.15'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.77'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.94%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.84'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.68%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +7.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0908'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.58%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.93'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.04'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.226.98'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.48%  '
